$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1635
$ws.Range("E2").Value = -1
$ws.Range("F2").Value = -1
$ws.Range("G2").Value = -29
$ws.Range("H2").Value = -35
$ws.Range("I2").Value = -35
$ws.Range("J2").ClearContents()
$ws.Range("K2").Value = 3346
$ws.Range("L2").Value = 849
$ws.Range("M2").Value = 2497
$ws.Range("N2").Value = 2497
$ws.Range("O2").ClearContents()
$ws.Range("P2").Value = 300
$ws.Range("Q2").Value = -66
$ws.Range("R2").Value = 94
$ws.Range("S2").Value = -26
$ws.Range("T2").Value = 117
$ws.Range("U2").Value = -183
$ws.Range("V2").Value = 417
$ws.Range("W2").Value = -0.06
$ws.Range("X2").Value = -2.16
$ws.Range("Y2").Value = -1.4
$ws.Range("Z2").Value = -1.07
$ws.Range("AA2").Value = 33.99
$ws.Range("AB2").Value = 732.4400000000001
$ws.Range("AC2").Value = -60
$ws.Range("AD2").Value = -31.82
$ws.Range("AE2").Value = 4243
$ws.Range("AF2").Value = 0.45
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 58855502

# Row 3
$ws.Range("D3").Value = 1773
$ws.Range("E3").Value = 31
$ws.Range("F3").Value = 31
$ws.Range("G3").Value = 20
$ws.Range("H3").Value = 22
$ws.Range("I3").Value = 22
$ws.Range("J3").ClearContents()
$ws.Range("K3").Value = 3625
$ws.Range("L3").Value = 836
$ws.Range("M3").Value = 2789
$ws.Range("N3").Value = 2789
$ws.Range("O3").ClearContents()
$ws.Range("P3").Value = 360
$ws.Range("Q3").Value = 139
$ws.Range("R3").Value = -132
$ws.Range("S3").Value = 238
$ws.Range("T3").Value = 36
$ws.Range("U3").Value = 103
$ws.Range("V3").Value = 364
$ws.Range("W3").Value = 1.73
$ws.Range("X3").Value = 1.25
$ws.Range("Y3").Value = 0.84
$ws.Range("Z3").Value = 0.64
$ws.Range("AA3").Value = 29.95
$ws.Range("AB3").Value = 677.51
$ws.Range("AC3").Value = 38
$ws.Range("AD3").Value = 84.97
$ws.Range("AE3").Value = 4028
$ws.Range("AF3").Value = 0.79
$ws.Range("AG3").Value = 30
$ws.Range("AH3").Value = 0.9399999999999999
$ws.Range("AI3").Value = 93.65000000000001
$ws.Range("AJ3").Value = 69751600

# Row 4
$ws.Range("D4").Value = 1814
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = 10
$ws.Range("G4").Value = 207
$ws.Range("H4").Value = 211
$ws.Range("I4").Value = 211
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 6364
$ws.Range("L4").Value = 1327
$ws.Range("M4").Value = 5038
$ws.Range("N4").Value = 5030
$ws.Range("O4").Value = 7
$ws.Range("P4").Value = 360
$ws.Range("Q4").Value = 23
$ws.Range("R4").Value = -190
$ws.Range("S4").Value = -67
$ws.Range("T4").Value = 94
$ws.Range("U4").Value = -71
$ws.Range("V4").Value = 359
$ws.Range("W4").Value = 0.54
$ws.Range("X4").Value = 11.64
$ws.Range("Y4").Value = 5.41
$ws.Range("Z4").Value = 4.23
$ws.Range("AA4").Value = 26.34
$ws.Range("AB4").Value = 730.33
$ws.Range("AC4").Value = 303
$ws.Range("AD4").Value = 10.49
$ws.Range("AE4").Value = 7313
$ws.Range("AF4").Value = 0.43
$ws.Range("AG4").Value = 10
$ws.Range("AH4").Value = 0.31
$ws.Range("AI4").Value = 3.25
$ws.Range("AJ4").Value = 69751600

# Row 5
$ws.Range("D5").Value = 1783
$ws.Range("E5").Value = -2
$ws.Range("F5").Value = -2
$ws.Range("G5").Value = 42
$ws.Range("H5").Value = 41
$ws.Range("I5").Value = 42
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 6541
$ws.Range("L5").Value = 1521
$ws.Range("M5").Value = 5020
$ws.Range("N5").Value = 5017
$ws.Range("O5").Value = 3
$ws.Range("P5").Value = 360
$ws.Range("Q5").Value = 67
$ws.Range("R5").Value = 124
$ws.Range("S5").Value = -181
$ws.Range("T5").Value = 286
$ws.Range("U5").Value = -219
$ws.Range("V5").Value = 621
$ws.Range("W5").Value = -0.13
$ws.Range("X5").Value = 2.31
$ws.Range("Y5").Value = 0.83
$ws.Range("Z5").Value = 0.64
$ws.Range("AA5").Value = 30.31
$ws.Range("AB5").Value = 741.64
$ws.Range("AC5").Value = 60
$ws.Range("AD5").Value = 43.98
$ws.Range("AE5").Value = 7498
$ws.Range("AF5").Value = 0.35
$ws.Range("AG5").ClearContents()
$ws.Range("AH5").ClearContents()
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 69751600

# Row 6
$ws.Range("D6").Value = 1611
$ws.Range("E6").Value = -86
$ws.Range("F6").Value = -86
$ws.Range("G6").Value = -92
$ws.Range("H6").Value = -79
$ws.Range("I6").Value = -79
$ws.Range("K6").Value = 6487
$ws.Range("L6").Value = 1593
$ws.Range("M6").Value = 4894
$ws.Range("N6").Value = 4892
$ws.Range("P6").Value = 360
$ws.Range("Q6").Value = -1
$ws.Range("R6").Value = -169
$ws.Range("S6").Value = 120
$ws.Range("T6").Value = 155
$ws.Range("U6").Value = -157
$ws.Range("V6").Value = 771
$ws.Range("W6").Value = -5.33
$ws.Range("X6").Value = -4.93
$ws.Range("Y6").Value = -1.6
$ws.Range("Z6").Value = -1.22
$ws.Range("AA6").Value = 32.55
$ws.Range("AB6").Value = 717.87
$ws.Range("AC6").Value = -113
$ws.Range("AD6").Value = -17.06
$ws.Range("AE6").Value = 7449
$ws.Range("AF6").Value = 0.26
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 69751600

# Row 7: clear all data cells D:AI
$ws.Range("D7:AI7").ClearContents()

# Row 8: clear all data cells D:AI
$ws.Range("D8:AI8").ClearContents()

# Row 9: clear all data cells D:AI
$ws.Range("D9:AI9").ClearContents()
